$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters that are updated, in order, with their new values for rows 2 & 3.
$values = @{
    "D"  = -0.602
    "G"  = 33.93939393939394
    "H"  = 33.93939393939394
    "I"  = -0.372960372960373
    "J"  = -0.372960372960373
    "K"  = -5.19
    "L"  = -1.20979020979021
    "M"  = 92.94
    "N"  = 0.5830614805520702
    "O"  = -17.90751445086705
    "P"  = 82.44
    "Q"  = 0.5171894604767879
    "R"  = -15.88439306358381
    "S"  = 10.5
    "T"  = 0.1129761136216914
    "U"  = 27
    "V"  = 0.1693851944792974
    "W"  = -0.02087691069991955
    "X"  = 0.04821709766681308
    "Y"  = -0.06909400836673263
    "Z"  = 0.02879561823319753
    "AA" = -0.01073962451587787
    "AB" = 0.04801944750369452
    "AC" = -0.05875907201957239
    "AD" = 2.08
    "AF" = 2.08
    "AG" = -24.92
    "AH" = 0.01288085211790934
    "AI" = 0.01406545847984852
    "AJ" = -0.1853063652587745
    "AK" = -0.2061548643282594
    "AL" = 0.147
    "AM" = 0.147
    "AN" = -1.71900826446281
    "AO" = -10.8843537414966
    "AP" = 20.59504132231405
    "AQ" = -10.8843537414966
}

foreach ($row in 2, 3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
    # Column E is removed entirely (was historical_growth_net_income_last_5_years value)
    $ws.Range("E$row").ClearContents()
}
